# Edit script applying the commit:
#   "Ajuste de la Queratina"  -> adjust keratin product commission split (rows 10,15,19,23,29,36)
#   "% a Marinela"            -> new Porc_trans/Cost_trans/Porc_producto/Valor_Neto columns
#   "Menor % a los profesionales de Tocador" -> recalculated Part_profesional values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 3 new columns at F (Porc_trans, Cost_trans, Porc_producto) ---
# old F (Valor_producto) / G (Part_profesional) / H (Revisar) shift right to I / J.. / K / L
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).Insert()

# --- Step 2: insert 1 new column at J (Valor_Neto), between Valor_producto(I) and Part_profesional ---
$ws.Columns.Item(10).Insert()

# --- Step 3: insert a brand-new data row at 29 (Sandra Giraldo / Blower cabello largo) ---
# this pushes the old "Fondo/Descuento" adjustment rows (29-34) down to (30-35)
$ws.Rows.Item(29).Insert()

# --- Step 4: header row - label the 4 new columns ---
$ws.Cells.Item(1,6).Value = "Porc_trans"
$ws.Cells.Item(1,7).Value = "Cost_trans"
$ws.Cells.Item(1,8).Value = "Porc_producto"
$ws.Cells.Item(1,10).Value = "Valor_Neto"

# --- Step 5: populate Porc_trans(F), Cost_trans(G), Porc_producto(H) and Valor_Neto(J) for every service/product row (2-29) ---
$ws.Cells.Item(2,6).Value = 0.036; $ws.Cells.Item(2,7).Value = 3240; $ws.Cells.Item(2,8).Value = 0.04534444444444444; $ws.Cells.Item(2,10).Value = 82679
$ws.Cells.Item(3,6).Value = 0.036; $ws.Cells.Item(3,7).Value = 1152; $ws.Cells.Item(3,8).Value = 0.26; $ws.Cells.Item(3,10).Value = 22528
$ws.Cells.Item(4,6).Value = 0.036; $ws.Cells.Item(4,7).Value = 791.9999999999999; $ws.Cells.Item(4,8).Value = 0.25; $ws.Cells.Item(4,10).Value = 15708
$ws.Cells.Item(5,6).Value = 0.036; $ws.Cells.Item(5,7).Value = 720; $ws.Cells.Item(5,8).Value = 0.20405; $ws.Cells.Item(5,10).Value = 15199
$ws.Cells.Item(6,6).Value = 0.036; $ws.Cells.Item(6,7).Value = 1260; $ws.Cells.Item(6,8).Value = 0.1166; $ws.Cells.Item(6,10).Value = 29659
$ws.Cells.Item(7,6).Value = 0.036; $ws.Cells.Item(7,7).Value = 1260; $ws.Cells.Item(7,8).Value = 0.1166; $ws.Cells.Item(7,10).Value = 29659
$ws.Cells.Item(8,6).Value = 0.036; $ws.Cells.Item(8,7).Value = 1980; $ws.Cells.Item(8,8).Value = 0.0742; $ws.Cells.Item(8,10).Value = 48939
$ws.Cells.Item(9,6).Value = 0.036; $ws.Cells.Item(9,7).Value = 6839.999999999999; $ws.Cells.Item(9,8).Value = 0.02147894736842105; $ws.Cells.Item(9,10).Value = 179079
$ws.Cells.Item(10,6).Value = 0.036; $ws.Cells.Item(10,7).Value = 5220; $ws.Cells.Item(10,8).Value = 0.3364413793103448; $ws.Cells.Item(10,10).Value = 90996
$ws.Cells.Item(11,6).Value = 0.036; $ws.Cells.Item(11,7).Value = 1620; $ws.Cells.Item(11,8).Value = 0.09068888888888889; $ws.Cells.Item(11,10).Value = 39299
$ws.Cells.Item(12,6).Value = 0; $ws.Cells.Item(12,7).Value = 0; $ws.Cells.Item(12,8).Value = 0.20405; $ws.Cells.Item(12,10).Value = 15919
$ws.Cells.Item(13,6).Value = 0.036; $ws.Cells.Item(13,7).Value = 1260; $ws.Cells.Item(13,8).Value = 0.1166; $ws.Cells.Item(13,10).Value = 29659
$ws.Cells.Item(14,6).Value = 0.036; $ws.Cells.Item(14,7).Value = 3240; $ws.Cells.Item(14,8).Value = 0.04534444444444444; $ws.Cells.Item(14,10).Value = 82679
$ws.Cells.Item(15,6).Value = 0.036; $ws.Cells.Item(15,7).Value = 6839.999999999999; $ws.Cells.Item(15,8).Value = 0.1610736842105263; $ws.Cells.Item(15,10).Value = 152556
$ws.Cells.Item(16,6).Value = 0.036; $ws.Cells.Item(16,7).Value = 1620; $ws.Cells.Item(16,8).Value = 0.09068888888888889; $ws.Cells.Item(16,10).Value = 39299
$ws.Cells.Item(17,6).Value = 0.036; $ws.Cells.Item(17,7).Value = 3600; $ws.Cells.Item(17,8).Value = 0.15; $ws.Cells.Item(17,10).Value = 81400
$ws.Cells.Item(18,6).Value = 0; $ws.Cells.Item(18,7).Value = 0; $ws.Cells.Item(18,8).Value = 0.1166; $ws.Cells.Item(18,10).Value = 30919
$ws.Cells.Item(19,6).Value = 0; $ws.Cells.Item(19,7).Value = 0; $ws.Cells.Item(19,8).Value = 0.2719242424242424; $ws.Cells.Item(19,10).Value = 96106
$ws.Cells.Item(20,6).Value = 0; $ws.Cells.Item(20,7).Value = 0; $ws.Cells.Item(20,8).Value = 0.1166; $ws.Cells.Item(20,10).Value = 30919
$ws.Cells.Item(21,6).Value = 0.036; $ws.Cells.Item(21,7).Value = 1386; $ws.Cells.Item(21,8).Value = 0.106; $ws.Cells.Item(21,10).Value = 33033
$ws.Cells.Item(22,6).Value = 0.036; $ws.Cells.Item(22,7).Value = 1620; $ws.Cells.Item(22,8).Value = 0.09068888888888889; $ws.Cells.Item(22,10).Value = 39299
$ws.Cells.Item(23,6).Value = 0.036; $ws.Cells.Item(23,7).Value = 14040; $ws.Cells.Item(23,8).Value = 0.1026666666666667; $ws.Cells.Item(23,10).Value = 335920
$ws.Cells.Item(24,6).Value = 0; $ws.Cells.Item(24,7).Value = 0; $ws.Cells.Item(24,8).Value = 0.09068888888888889; $ws.Cells.Item(24,10).Value = 40919
$ws.Cells.Item(25,6).Value = 0.036; $ws.Cells.Item(25,7).Value = 1386; $ws.Cells.Item(25,8).Value = 0.106; $ws.Cells.Item(25,10).Value = 33033
$ws.Cells.Item(26,6).Value = 0.036; $ws.Cells.Item(26,7).Value = 1260; $ws.Cells.Item(26,8).Value = 0.1166; $ws.Cells.Item(26,10).Value = 29659
$ws.Cells.Item(27,6).Value = 0.036; $ws.Cells.Item(27,7).Value = 3240; $ws.Cells.Item(27,8).Value = 0.04534444444444444; $ws.Cells.Item(27,10).Value = 82679
$ws.Cells.Item(28,6).Value = 0; $ws.Cells.Item(28,7).Value = 0; $ws.Cells.Item(28,8).Value = 0.106; $ws.Cells.Item(28,10).Value = 34419
$ws.Cells.Item(29,6).Value = 0; $ws.Cells.Item(29,7).Value = 0; $ws.Cells.Item(29,8).Value = 0.09068888888888889; $ws.Cells.Item(29,10).Value = 40919

# --- Step 6: rows where a keratin product charge now reduces Part_profesional (K) ---
# Valor_producto (I) and Part_profesional (K) both change on these 4 rows
$ws.Cells.Item(10,9).Value = 48784; $ws.Cells.Item(10,11).Value = 30966
$ws.Cells.Item(15,9).Value = 30604; $ws.Cells.Item(15,11).Value = 73896
$ws.Cells.Item(19,9).Value = 35894; $ws.Cells.Item(19,11).Value = 36706
$ws.Cells.Item(23,9).Value = 40040; $ws.Cells.Item(23,11).Value = 174460

# --- Step 7: fill in the brand-new row 29 (Sandra Giraldo / Blower cabello largo) ---
$ws.Cells.Item(29,1).Value = "01/11/2024 11:43"
$ws.Cells.Item(29,2).Value = "Sandra Giraldo"
$ws.Cells.Item(29,3).Value = "Blower cabello largo"
$ws.Cells.Item(29,4).Value = "Olga Arango"
$ws.Cells.Item(29,5).Value = 45000
$ws.Cells.Item(29,6).Value = 0
$ws.Cells.Item(29,7).Value = 0
$ws.Cells.Item(29,8).Value = 0.09068888888888889
$ws.Cells.Item(29,9).Value = 4081
$ws.Cells.Item(29,10).Value = 40919
$ws.Cells.Item(29,11).Value = 24750

# --- Step 8: append brand-new row 36 (Descuento - Producto - Productos de Color) ---
# leading apostrophe forces literal text (this date-like string would otherwise auto-convert to a date)
$ws.Cells.Item(36,1).Value = "'2024-11-01"
$ws.Cells.Item(36,3).Value = "Descuento - Producto - Productos de Color"
$ws.Cells.Item(36,4).Value = "Olga Arango"
$ws.Cells.Item(36,11).Value = -58905

Write-Host "edit complete"
